$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Range("Z4").Interior.ThemeColor = 2
$ws.Range("Z4").Interior.TintAndShade = -0.0999786370433668
